$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections on row 13 (Nature Biomedical Engineering *) ---
$ws.Range("G13").Value = 21
$ws.Range("H13").Value = 22
$ws.Range("I13").Value = 4

# --- K3: re-enter as its own (non-shared) formula ---
$ws.Range("K3").Formula = "=D3/B3*100"

# --- Totals row (row 15): sums now cover rows 2:14 instead of 5:14 ---
$ws.Range("B15").Formula = "=SUM(B2:B14)"
$ws.Range("C15").Formula = "=SUM(C2:C14)"
$ws.Range("E15").Formula = "=SUM(E2:E14)"
$ws.Range("F15").Formula = "=SUM(F2:F14)"
$ws.Range("G15").Formula = "=SUM(G2:G14)"
$ws.Range("H15").Formula = "=SUM(H2:H14)"
$ws.Range("I15").Formula = "=SUM(I2:I14)"

# --- View state: zoom + active selection ---
$excel.ActiveWindow.Zoom = 150
[void]$ws.Range("J13").Select()
